$wb = $excel.ActiveWorkbook

$ws1 = $wb.Sheets.Item("NewUsers")

# Insert the new "marketingStatus" row between the existing "graduationDate"
# row (9) and the "Token" section (formerly row 11, now shifted to 12).
$ws1.Rows.Item(10).Insert()

$ws1.Range("A10").Value = "marketingStatus"

$ws1.Range("B10").Style = "Normal"
$ws1.Range("B10").Value = $true

$ws1.Range("C10").Style = "Normal"
$ws1.Range("C10").Value = $true

$ws1.Range("D10").Value = "STOP"

# NewUsers becomes the active sheet / tab, with a new selection.
$ws1.Activate() | Out-Null
$ws1.Range("E17").Select() | Out-Null
